$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8136539459228516
$ws.Range("B1").Value = 2.867855787277222
$ws.Range("C1").Value = 7.588922500610352
$ws.Range("D1").Value = 2.203446865081787
$ws.Range("E1").Value = 1.473876476287842
